# Balancing Sheet update: add Damage Type per attack (Primary/Secondary) and
# change several values according to the updated design spreadsheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Balancing Sheet (2)")
$ws.Activate()

# --- Top stat block -------------------------------------------------------
# Stagger (row 4): Enemy value 100 -> 150
$ws.Range("C4").Value = 150

# Stamina Regen /s (row 7): Enemy value 10 -> 5
$ws.Range("C7").Value = 5

# --- DMG / Stamina Cost / Stagger Multiplier table (rows 14-19) ----------
# Primary Sword: Stagger Multiplier 1.2 -> 1.1
$ws.Range("D14").Value = 1.1

# Primary Hammer: DMG 90 -> 75 (now a plain number, left aligned)
$ws.Range("B16").Value = 75
$ws.Range("B16").HorizontalAlignment = -4131  # xlLeft

# Primary Hammer: Stagger Multiplier 1.8 -> 1.5
$ws.Range("D16").Value = 1.5

# Secondary Hammer: Stagger Multiplier 1.8 -> 2
$ws.Range("D17").Value = 2

# Primary Bow: Stagger Multiplier 1 -> 0.5
$ws.Range("D18").Value = 0.5

# Secondary Bow: Stagger Multiplier 1 -> 0.7
$ws.Range("D19").Value = 0.7

# --- View state: zoom + selection ----------------------------------------
$ws.Range("D19").Select()
$excel.ActiveWindow.Zoom = 85

$wb.Save()
